$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Foundation class: remove the old single "List<Card>" Type value from
# the Property row (row 21) and move it under the Property column of the
# following row (22), then add new Foundation properties (IsEmpty/bool,
# "FoundationPile (4x)"/List<Card>, OpenCards/List<Card>, ClosedCards/List<Card>).
# --- Tableau class: similarly shift its single Type value down one row.

# Clear out everything from row 19 through the end of the old "Board" class
# (row 28) -- contents AND formatting -- so we can re-lay the content out
# with the extra Foundation rows inserted and the Board class moved down to
# make room.
$ws.Range("A19:B28").Clear()

# Tableau class (rows 18-19)
$ws.Range("A19").Value = "TableauCards"
$ws.Range("B19").Value = "List<Card>"

# Foundation class (rows 21-26)
$ws.Range("A21").Value = "Foundation"
$ws.Range("A22").Value = "FoundationCards"
$ws.Range("B22").Value = "List<Card>"
$ws.Range("A23").Value = "IsEmpty"
$ws.Range("B23").Value = "bool"
$ws.Range("A24").Value = "FoundationPile (4x)"
$ws.Range("B24").Value = "List<Card>"
$ws.Range("A25").Value = "OpenCards"
$ws.Range("B25").Value = "List<Card>"
$ws.Range("A26").Value = "ClosedCards"
$ws.Range("B26").Value = "List<Card>"

# Board class (rows 28-32), moved down from rows 24-28
$ws.Range("A28").Value = "Board"
$ws.Range("A29").Value = "Tableau"
$ws.Range("B29").Value = "Tableau"
$ws.Range("A30").Value = "Stock"
$ws.Range("B30").Value = "Stock"
$ws.Range("A31").Value = "Waste"
$ws.Range("B31").Value = "Waste"
$ws.Range("A32").Value = "Foundation"
$ws.Range("B32").Value = "Foundation"

# Re-apply the "class header" shading style (style used by A4/A8/A12/A18/A21)
# to the relocated/row-shifted class header cells A21 and A28.
$ws.Range("A18").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A18").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Clear old header styling left behind on the now plain rows.
$ws.Range("B18").ClearContents()

# Column A needs to be a bit wider to fit the new, longer "FoundationPile (4x)" label.
$ws.Columns("A:A").ColumnWidth = 18.7109375

# Scroll / selection state, matching the author's saved view.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A23").Select()
